# Daily attendance processing - 2026-01-24 21:56:48
#
# The "Recorded By" column (G) lists the users who recorded/edited a
# session's attendance, as a comma-separated string. This pass rotates
# the list for every multi-author cell so that the most recent editor
# (previously listed last) is moved to the front of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Text

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    if ($val.IndexOf(",") -lt 0) {
        continue
    }

    $rawParts = $val.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    if ($parts.Count -lt 2) {
        continue
    }

    $lastEntry = $parts[$parts.Count - 1]
    $remaining = $parts[0..($parts.Count - 2)]
    $rotated = @($lastEntry) + $remaining
    $newVal = $rotated -join ", "

    $cell.Value = $newVal
}
